$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "Marengoni e Stringhini (2009)" (row 3).
# All rows below shift up by one.
$ws.Rows.Item(3).Delete()

# After the delete, the table rows (Autores) are now:
#   3 Pereira (2017)
#   4 Bispo (2018)
#   5 Fernandes (2019)
#   6 Santos (2019)
#   7 Fonseca (2020)
#   8 Cruz (2020)
#   9 Nascimento et al. (2020)
#  10 Nogueira (2020)
#  11 Júnior (2021)
#  12 Martins (2021)

# Fill in the new "Trabalhos Relacionados" entries about LSTM-based works.
# Order of entry matters for shared-string creation order.

# Bispo (2018) - row 4
$ws.Range("C4").Value2 = "LSTM;"
$ws.Range("D4").Value2 = "Criou um modelo para identificar comentários ofensivos na língua portuguesa, baseando-se em um modelo já existente porém com dados da língua inglesa;"

# Santos (2019) - row 6
$ws.Range("C6").Value2 = "LSTM;"
$ws.Range("D6").Value2 = "Realizou um trabalho de predição de preços de energia elétrica utilizando um modelo LSTM, comparando os resultados não só com os dados reais, mas também com os resultados de um modelo simples em Bootstrap;"

# Júnior (2021) - row 11
$ws.Range("C11").Value2 = "LSTM;"
$ws.Range("D11").Value2 = "Desenvolveu uma ferramenta web para predição no mercado de ações, utilizando uma rede LSTM como base em sua execução;`nDemonstra durante o desenvolvimento de seu trabalho as dificuldades encontradas e como elas foram mitigadas;"
$ws.Range("E11").Value2 = "Não utilizou nenhum outro modelo para realizar comparativos de resultados. Com isso, não fica confirmada completamente a eficácia do modelo LSTM implementado;"

# Pereira (2017) - row 3
$ws.Range("C3").Value2 = "LSTM;"
$ws.Range("D3").Value2 = "Realizou uma análise comparativa de resultados dos modelos LSTM, ARIMA e RNA, exibindo a efetividade do modelo LSTM para a situação aplicada;"

# Fernandes (2019) - row 5 (only the topic column is filled)
$ws.Range("C5").Value2 = "LSTM;"

# Adjust row heights to match the new content layout.
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 45.75

# Adjust column widths (B narrower, D/E much wider to fit the new long texts).
$ws.Columns.Item(2).ColumnWidth = 21.75
$ws.Columns.Item(4).ColumnWidth = 77.6
$ws.Columns.Item(5).ColumnWidth = 76.5

# Update the view: zoom to 100% and move the selection.
$excel.ActiveWindow.Zoom = 100
$ws.Range("K12").Select() | Out-Null
